# Add two new localization rows to Sheet1:
#   row 91: lang_parent   | Phụ Huynh | Parent
#   row 92: lang_student  | Học Sinh  | Student

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A91").Value = "lang_parent"
$ws.Range("B91").Value = "Phụ Huynh"
$ws.Range("C91").Value = "Parent"

$ws.Range("A92").Value = "lang_student"
$ws.Range("B92").Value = "Học Sinh"
$ws.Range("C92").Value = "Student"

# Apply the same formatting as the rest of the data rows (cellXfs style
# index 3, i.e. default font with wrap-text alignment) to the new cells.
$ws.Range("A91:C92").WrapText = $true

# Match the recorded view state from the edit: scroll position stays,
# but the active selection moves to C92.
$ws.Range("C92").Select()
